$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.143.88"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.639.70"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'0.992"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").Value = "'215.88"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "'0.0635"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'19.84"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'4.26"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.866.02"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "1.626.50"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'0.554"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'63.42"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "26.105.56"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D20").Value = "'4.47"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'193.92"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'10.01"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'6.38"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'141.90"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'6.91"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'15.62"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'0.0497"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'1.60"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'0.910"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "1.145.34"
$ws.Range("D38").Value = "'0.548"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'0.0157"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'5.57"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'100.38"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.795"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.774.59"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "'55.91"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "0.0₆0106"
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "'1.46"
$ws.Range("E49").Value = "  +5.01%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0974"
$ws.Range("E51").Value = "  +2.94%  "
